# Corrected some selection scopes
# Remove every other data row (rows 2,4,6,...,30) from Sheet1 so the
# remaining rows shift up and get renumbered sequentially.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsToDelete = @(2,4,6,8,10,12,14,16,18,20,22,24,26,28,30)

# Delete from the bottom up so earlier row numbers stay valid as we go.
$sorted = $rowsToDelete | Sort-Object -Descending
foreach ($r in $sorted) {
    $ws.Rows.Item($r).Delete()
}
